# "fix call stack and double quote character"
#
# 1) Handout/Notes master date footer: refresh the cached "last saved" date
#    text (vendredi 17 janvier 2025 -> mardi 21 janvier 2025).
# 2) Slide 1 "Rectangle 2" code sample: replace the stray left smart-quote
#    before c_code with a straight double-quote ("c_code" -> "c_code").
# 3) Slide 1 "Rectangle 4" (the call-stack box outline): narrow it so it no
#    longer overlaps the neighbouring diagram.
# 4) Slide 2 "TextBox 13": rename method2 -> function2 to match the sibling
#    function1 textbox used in the call-stack diagram.

$p = $ppt.ActivePresentation

# --- 1. Date footers on Handout Master and Notes Master -------------------
$handoutMaster = $p.HandoutMaster
$handoutMaster.HeadersFooters.DateAndTime.Text = "mardi 21 janvier 2025"

$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = "mardi 21 janvier 2025"

# --- 2. Smart-quote fix on slide 1 -----------------------------------------
$slide1 = $p.Slides.Item(1)
$codeShape = $slide1.Shapes.Item(1)
$codeText = $codeShape.TextFrame.TextRange
$quoteIdx = $codeText.Text.IndexOf('"c_code"') + 1
$codeText.Characters($quoteIdx, 1).Text = [string][char]0x22

# --- 3. Resize the call-stack rectangle on slide 1 -------------------------
$callStackShape = $slide1.Shapes.Item(2)
$callStackShape.Left = 461.52058
$callStackShape.Width = 147.7624

# --- 4. Rename method2 -> function2 on slide 2 -----------------------------
$slide2 = $p.Slides.Item(2)
$functionShape = $slide2.Shapes.Item(12)
$functionShape.TextFrame.TextRange.Text = "function2"
